$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6035
$ws.Range("E2").Value = 185
$ws.Range("F2").Value = 185
$ws.Range("G2").Value = 166
$ws.Range("H2").Value = 520
$ws.Range("I2").Value = 520
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4833
$ws.Range("L2").Value = 2059
$ws.Range("M2").Value = 2774
$ws.Range("N2").Value = 2774
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 103
$ws.Range("Q2").Value = 414
$ws.Range("R2").Value = -235
$ws.Range("S2").Value = 13
$ws.Range("T2").Value = 356
$ws.Range("U2").Value = 58
$ws.Range("V2").Value = 1069
$ws.Range("W2").Value = 3.06
$ws.Range("X2").Value = 8.619999999999999
$ws.Range("Y2").Value = 16.85
$ws.Range("Z2").Value = 9.460000000000001
$ws.Range("AA2").Value = 74.23
$ws.Range("AB2").Value = 3657.16
$ws.Range("AC2").Value = 1689
$ws.Range("AD2").Value = 4.59
$ws.Range("AE2").Value = 14090
$ws.Range("AF2").Value = 0.55
$ws.Range("AG2").Value = 230
$ws.Range("AH2").Value = 2.96
$ws.Range("AI2").Value = 8.699999999999999
$ws.Range("AJ2").Value = 20623600

# Row 3
$ws.Range("D3").Value = 6424
$ws.Range("E3").Value = 252
$ws.Range("F3").Value = 252
$ws.Range("G3").Value = 207
$ws.Range("H3").Value = 162
$ws.Range("I3").Value = 178
$ws.Range("J3").Value = -15
$ws.Range("K3").Value = 7004
$ws.Range("L3").Value = 3082
$ws.Range("M3").Value = 3922
$ws.Range("N3").Value = 3242
$ws.Range("O3").Value = 681
$ws.Range("P3").Value = 132
$ws.Range("Q3").Value = 449
$ws.Range("R3").Value = -382
$ws.Range("S3").Value = 108
$ws.Range("T3").Value = 451
$ws.Range("U3").Value = -2
$ws.Range("V3").Value = 1922
$ws.Range("W3").Value = 3.92
$ws.Range("X3").Value = 2.53
$ws.Range("Y3").Value = 5.92
$ws.Range("Z3").Value = 2.74
$ws.Range("AA3").Value = 78.56999999999999
$ws.Range("AB3").Value = 2393.26
$ws.Range("AC3").Value = 814
$ws.Range("AD3").Value = 8.6
$ws.Range("AE3").Value = 12772
$ws.Range("AF3").Value = 0.55
$ws.Range("AG3").Value = 110
$ws.Range("AH3").Value = 1.57
$ws.Range("AI3").Value = 15.69
$ws.Range("AJ3").Value = 26319633

# Row 4
$ws.Range("D4").Value = 7234
$ws.Range("E4").Value = 363
$ws.Range("F4").Value = 363
$ws.Range("G4").Value = 261
$ws.Range("H4").Value = 220
$ws.Range("I4").Value = 149
$ws.Range("J4").Value = 71
$ws.Range("K4").Value = 7224
$ws.Range("L4").Value = 3183
$ws.Range("M4").Value = 4041
$ws.Range("N4").Value = 3308
$ws.Range("O4").Value = 733
$ws.Range("P4").Value = 132
$ws.Range("Q4").Value = 621
$ws.Range("R4").Value = -370
$ws.Range("S4").Value = -102
$ws.Range("T4").Value = 443
$ws.Range("U4").Value = 178
$ws.Range("V4").Value = 1870
$ws.Range("W4").Value = 5.02
$ws.Range("X4").Value = 3.03
$ws.Range("Y4").Value = 4.54
$ws.Range("Z4").Value = 3.08
$ws.Range("AA4").Value = 78.77
$ws.Range("AB4").Value = 2480.82
$ws.Range("AC4").Value = 565
$ws.Range("AD4").Value = 12.54
$ws.Range("AE4").Value = 13033
$ws.Range("AF4").Value = 0.54
$ws.Range("AG4").Value = 140
$ws.Range("AH4").Value = 1.97
$ws.Range("AI4").Value = 23.88
$ws.Range("AJ4").Value = 26319633

# Row 5
$ws.Range("D5").Value = 8223
$ws.Range("E5").Value = 454
$ws.Range("F5").Value = 454
$ws.Range("G5").Value = 452
$ws.Range("H5").Value = 409
$ws.Range("I5").Value = 208
$ws.Range("J5").Value = 201
$ws.Range("K5").Value = 7399
$ws.Range("L5").Value = 3074
$ws.Range("M5").Value = 4325
$ws.Range("N5").Value = 3410
$ws.Range("O5").Value = 916
$ws.Range("P5").Value = 132
$ws.Range("Q5").Value = 250
$ws.Range("R5").Value = -388
$ws.Range("S5").Value = -120
$ws.Range("T5").Value = 412
$ws.Range("U5").Value = -161
$ws.Range("V5").Value = 1726
$ws.Range("W5").Value = 5.52
$ws.Range("X5").Value = 4.97
$ws.Range("Y5").Value = 6.19
$ws.Range("Z5").Value = 5.59
$ws.Range("AA5").Value = 71.06999999999999
$ws.Range("AB5").Value = 2616.97
$ws.Range("AC5").Value = 790
$ws.Range("AD5").Value = 8.140000000000001
$ws.Range("AE5").Value = 13433
$ws.Range("AF5").Value = 0.48
$ws.Range("AG5").Value = 120
$ws.Range("AH5").Value = 1.87
$ws.Range("AI5").Value = 14.65
$ws.Range("AJ5").Value = 26319633

# Row 6
$ws.Range("D6").Value = 8910
$ws.Range("E6").Value = 463
$ws.Range("F6").Value = 463
$ws.Range("G6").Value = 409
$ws.Range("H6").Value = 227
$ws.Range("I6").Value = 60
$ws.Range("K6").Value = 7425
$ws.Range("L6").Value = 2967
$ws.Range("M6").Value = 4458
$ws.Range("N6").Value = 3397
$ws.Range("P6").Value = 132
$ws.Range("Q6").Value = 532
$ws.Range("R6").Value = -607
$ws.Range("S6").Value = -334
$ws.Range("T6").Value = 373
$ws.Range("U6").Value = 159
$ws.Range("V6").Value = 1462
$ws.Range("W6").Value = 5.2
$ws.Range("X6").Value = 2.54
$ws.Range("Y6").Value = 1.75
$ws.Range("Z6").Value = 3.06
$ws.Range("AA6").Value = 66.55
$ws.Range("AB6").Value = 2617.91
$ws.Range("AC6").Value = 227
$ws.Range("AD6").Value = 24.13
$ws.Range("AE6").Value = 13385
$ws.Range("AF6").Value = 0.41
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 2.74
$ws.Range("AI6").Value = 63.82
$ws.Range("AJ6").Value = 26319633

# Remove the forecast rows 7-9 (2019E/2020E/2021E) financial data; keep only the row label columns A-C
$ws.Range("D7:AJ9").ClearContents()
